# Add 2023 (column S) data and a trailing blank column (U) to the
# "vachroba" regional trade worksheet, mirroring the formatting already
# used by the neighbouring columns R (last data year) and T (trailing
# blank column).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Copy formatting for the new column S (year 2023) from column R ---
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)

$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)

$ws.Range("R15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Copy formatting for the new trailing blank column U from column T ---
$ws.Range("T2:T15").Copy()
$ws.Range("U2:U15").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Populate the new 2023 values ---
$ws.Range("S3").Value = 2023
$ws.Range("S4").Value = 810.5
$ws.Range("S5").Value = 135.19999999999999
$ws.Range("S6").Value = 3146
$ws.Range("S7").Value = 2339
$ws.Range("S8").Value = 1503.1
$ws.Range("S9").Value = 41.9
$ws.Range("S10").Value = 42.8
$ws.Range("S11").Value = 93.3
$ws.Range("S12").Value = 5.9
$ws.Range("S13").Value = 746
$ws.Range("S14").Value = 704.7

# --- Extend the header merge to cover the new column ---
$ws.Range("A1:S1").Merge()

# --- Update the selected range shown in the sheet view ---
$ws.Range("S3:S14").Select()
